$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Update the "time_taken" (column F) timestamps on the "data" sheet ---
$newTimes = @(
    "2021-10-05 14:19:49.980652",
    "2021-10-05 14:19:49.980660",
    "2021-10-05 14:19:49.980663",
    "2021-10-05 14:19:49.980666",
    "2021-10-05 14:19:49.980668",
    "2021-10-05 14:19:49.980671",
    "2021-10-05 14:19:49.980674",
    "2021-10-05 14:19:49.980676",
    "2021-10-05 14:19:49.980679",
    "2021-10-05 14:19:49.980681",
    "2021-10-05 14:19:49.980684",
    "2021-10-05 14:19:49.980686",
    "2021-10-05 14:19:49.980689",
    "2021-10-05 14:19:49.980691",
    "2021-10-05 14:19:49.980694",
    "2021-10-05 14:19:49.980697",
    "2021-10-05 14:19:49.980700",
    "2021-10-05 14:19:49.980702",
    "2021-10-05 14:19:49.980705",
    "2021-10-05 14:19:49.980707",
    "2021-10-05 14:19:49.980710",
    "2021-10-05 14:19:49.980713",
    "2021-10-05 14:19:49.980715",
    "2021-10-05 14:19:49.980717",
    "2021-10-05 14:19:49.980720",
    "2021-10-05 14:19:49.980723",
    "2021-10-05 14:19:49.980726",
    "2021-10-05 14:19:49.980728",
    "2021-10-05 14:19:49.980731",
    "2021-10-05 14:19:49.980733",
    "2021-10-05 14:19:49.980736",
    "2021-10-05 14:19:49.980738",
    "2021-10-05 14:19:49.980741",
    "2021-10-05 14:19:49.980744",
    "2021-10-05 14:19:49.980746",
    "2021-10-05 14:19:49.980749",
    "2021-10-05 14:19:49.980752",
    "2021-10-05 14:19:49.980754",
    "2021-10-05 14:19:49.980757",
    "2021-10-05 14:19:49.980759",
    "2021-10-05 14:19:49.980762",
    "2021-10-05 14:19:49.980765",
    "2021-10-05 14:19:49.980768",
    "2021-10-05 14:19:49.980770",
    "2021-10-05 14:19:49.980773",
    "2021-10-05 14:19:49.980775",
    "2021-10-05 14:19:49.980778",
    "2021-10-05 14:19:49.980780",
    "2021-10-05 14:19:49.980783",
    "2021-10-05 14:19:49.980785",
    "2021-10-05 14:19:49.980788",
    "2021-10-05 14:19:49.980790",
    "2021-10-05 14:19:49.980793",
    "2021-10-05 14:19:49.980796",
    "2021-10-05 14:19:49.980798",
    "2021-10-05 14:19:49.980801",
    "2021-10-05 14:19:49.980803",
    "2021-10-05 14:19:49.980806",
    "2021-10-05 14:19:49.980808",
    "2021-10-05 14:19:49.980811",
    "2021-10-05 14:19:49.980813",
    "2021-10-05 14:19:49.980816",
    "2021-10-05 14:19:49.980818",
    "2021-10-05 14:19:49.980821",
    "2021-10-05 14:19:49.980825",
    "2021-10-05 14:19:49.980827",
    "2021-10-05 14:19:49.980830",
    "2021-10-05 14:19:49.980833",
    "2021-10-05 14:19:49.980835",
    "2021-10-05 14:19:49.980837",
    "2021-10-05 14:19:49.980840",
    "2021-10-05 14:19:49.980842",
    "2021-10-05 14:19:49.980845",
    "2021-10-05 14:19:49.980847",
    "2021-10-05 14:19:49.980850",
    "2021-10-05 14:19:49.980852",
    "2021-10-05 14:19:49.980857",
    "2021-10-05 14:19:49.980860",
    "2021-10-05 14:19:49.980862",
    "2021-10-05 14:19:49.980865",
    "2021-10-05 14:19:49.980867",
    "2021-10-05 14:19:49.980870",
    "2021-10-05 14:19:49.980872",
    "2021-10-05 14:19:49.980875",
    "2021-10-05 14:19:49.980878",
    "2021-10-05 14:19:49.980880",
    "2021-10-05 14:19:49.980883",
    "2021-10-05 14:19:49.980886",
    "2021-10-05 14:19:49.980889",
    "2021-10-05 14:19:49.980891",
    "2021-10-05 14:19:49.980894",
    "2021-10-05 14:19:49.980896",
    "2021-10-05 14:19:49.980900",
    "2021-10-05 14:19:49.980903",
    "2021-10-05 14:19:49.980906",
    "2021-10-05 14:19:49.980909",
    "2021-10-05 14:19:49.980911",
    "2021-10-05 14:19:49.980914",
    "2021-10-05 14:19:49.980917",
    "2021-10-05 14:19:49.980919",
    "2021-10-05 14:19:49.980922",
    "2021-10-05 14:19:49.980925",
    "2021-10-05 14:19:49.980927",
    "2021-10-05 14:19:49.980930",
    "2021-10-05 14:19:49.980933",
    "2021-10-05 14:19:49.980935",
    "2021-10-05 14:19:49.980938",
    "2021-10-05 14:19:49.980940",
    "2021-10-05 14:19:49.980945",
    "2021-10-05 14:19:49.980948",
    "2021-10-05 14:19:49.980951",
    "2021-10-05 14:19:49.980954",
    "2021-10-05 14:19:49.980956",
    "2021-10-05 14:19:49.980959",
    "2021-10-05 14:19:49.980962",
    "2021-10-05 14:19:49.980964",
    "2021-10-05 14:19:49.980967",
    "2021-10-05 14:19:49.980970",
    "2021-10-05 14:19:49.980973",
    "2021-10-05 14:19:49.980975",
    "2021-10-05 14:19:49.980978",
    "2021-10-05 14:19:49.980980",
    "2021-10-05 14:19:49.980983",
    "2021-10-05 14:19:49.980985",
    "2021-10-05 14:19:49.980988",
    "2021-10-05 14:19:49.980990",
    "2021-10-05 14:19:49.980993",
    "2021-10-05 14:19:49.980996",
    "2021-10-05 14:19:49.981000"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- 2. Add the new "metadata" sheet after "data" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadata"

# Reuse the bold/bordered header style already used on the "data" sheet's
# header row, and the index-column style used in column A, so the new
# sheet's styling matches the rest of the workbook.
$ws1.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Craniosynostosis"
$ws2.Range("C2").Value = 168

# "2.53" must stay textual (not be coerced to the number 2.53), so format
# the cell as Text before assigning the value.
$verCell = $ws2.Cells.Item(2, 4)
$verCell.NumberFormat = "@"
$verCell.Value = "2.53"

$ws2.Range("E2").Value = "2021-09-14T14:22:31.109410Z"
$ws2.Range("F2").Value = "2021-10-05 14:19:49.977376"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/168/?format=json"
